$d = $word.ActiveDocument

# The document previously carried a stray "_GoBack" bookmark right after
# "Пример корректной строки:" (a leftover from Word's last-edit-location
# tracking). The edit relocates it to mark the point the author actually
# last typed at: inside the first paragraph, right after "считывающ" in
# "...под считывающей головкой." Remove it from its old spot first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert the missing clause "Перед выполнением программы " before
# "установите", lower-casing the "У" so the sentence still reads
# naturally: "...на клавиатуре. Перед выполнением программы установите
# нужный символ под считывающей головкой."
$d.Content.Find.Execute("Установите нужный символ под считывающей головкой.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Перед выполнением программы установите нужный символ под считывающей головкой.", 2)

# Re-create "_GoBack" at the new location: right after "считывающ", i.e.
# immediately before "ей головкой."
$r = $d.Content
$r.Find.Execute("считывающ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $r.Duplicate
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
